$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 22482006
$ws.Range("L7").Value = 134036
$ws.Range("B10").Value = 22482007

$ws.Range("I16").Select()
